$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column I (index 9) values from 4 to 5 for rows 2-5
$ws.Cells.Item(2, 9).Value = 5
$ws.Cells.Item(3, 9).Value = 5
$ws.Cells.Item(4, 9).Value = 5
$ws.Cells.Item(5, 9).Value = 5

# Add new row 6 with data
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 12
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 10).Value = "train_dim2_1"

# Update the selected cell to I7 (as reflected in the diff)
$ws.Range("I7").Select()
